# Network diagram: add the "Applications -> App 3" dependency rows that are
# generated when the user clicks the App 3 node in the diagram's right-hand
# container (same shape as the existing App 1/App 2 dependency blocks), and
# leave the new row as the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("Applications", "App 3", "Parent Description…", "Depends On", "Facilities",    "Location 3", "Dependency Description…"),
    @("Applications", "App 3", "Parent Description…", "Depends On", "Procurements",  "PO 3",       "Dependency Description…"),
    @("Applications", "App 3", "Parent Description…", "Depends On", "People",        "Person 3",   "Dependency Description…"),
    @("Applications", "App 3", "Parent Description…", "Depends On", "Data",          "Data 3",     "Dependency Description…"),
    @("Applications", "App 3", "Parent Description…", "Depends On", "Technology",    "Tech 3",     "Dependency Description…")
)

$startRow = 37
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowValues = $newRows[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $rowValues[$j]
    }
}

# Scroll the view back to the top and make the freshly-clicked node's cell
# (the active node) the current selection.
$ws.Range("F44").Select()
